# "Generate Report for Archive"
#
# The localization-status report is regenerated: the shared "Ready for
# handoff" status text becomes "In Translation" everywhere it appears
# (the Overview roll-up sheet's per-locale status columns, plus each
# locale sheet's own Status column), and the now-shorter status text
# causes the affected status columns to shrink slightly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
